$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 19 (ALC)
$ws_ALC.Range("H19").Value = 2463.889
$ws_ALC.Range("I19").Value = 2106.4285
$ws_ALC.Range("J19").Value = 2691.3635
$ws_ALC.Range("K19").Value = 2106.4285
$ws_ALC.Range("L19").Value = 2691.3635
$ws_ALC.Range("M19").Value = -1931.4285
$ws_ALC.Range("N19").Value = -3041.3635

# Row 52 (ALC)
$ws_ALC.Range("H52").Value = 2520.7144
$ws_ALC.Range("I52").Value = 3022.3333
$ws_ALC.Range("J52").Value = 2144.5
$ws_ALC.Range("K52").Value = 9066.999899999999
$ws_ALC.Range("L52").Value = 6433.5
$ws_ALC.Range("M52").Value = -8906.999899999999
$ws_ALC.Range("N52").Value = -6753.5

# Row 70 (ALC)
$ws_ALC.Range("H70").Value = 127397.75
$ws_ALC.Range("J70").Value = 251985.75
$ws_ALC.Range("L70").Value = 755957.25
$ws_ALC.Range("N70").Value = -756497.25

# Row 73 (ALC)
$ws_ALC.Range("H73").Value = 127397.75
$ws_ALC.Range("J73").Value = 251985.75
$ws_ALC.Range("L73").Value = 755957.25
$ws_ALC.Range("N73").Value = -757829.25

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws_ARM.Range("H61").Value = 6940.5713
$ws_ARM.Range("I61").Value = 6940.5713
$ws_ARM.Range("K61").Value = 6940.5713
$ws_ARM.Range("M61").Value = -6728.5713

# Row 63 (ARM)
$ws_ARM.Range("H63").Value = 6488.087
$ws_ARM.Range("I63").Value = 5725.095
$ws_ARM.Range("K63").Value = 5725.095
$ws_ARM.Range("M63").Value = -5039.095

# Row 66 (ARM)
$ws_ARM.Range("H66").Value = 6488.087
$ws_ARM.Range("I66").Value = 5725.095
$ws_ARM.Range("K66").Value = 28625.475
$ws_ARM.Range("M66").Value = -25193.475

# Row 74 (ARM)
$ws_ARM.Range("H74").Value = 2904.8235
$ws_ARM.Range("I74").Value = 2018.1538
$ws_ARM.Range("K74").Value = 2018.1538
$ws_ARM.Range("M74").Value = -1144.1538

# Row 77 (ARM)
$ws_ARM.Range("H77").Value = 2904.8235
$ws_ARM.Range("I77").Value = 2018.1538
$ws_ARM.Range("K77").Value = 10090.769
$ws_ARM.Range("M77").Value = -5722.769

# Row 97 (ARM)
$ws_ARM.Range("H97").Value = 2962.4375
$ws_ARM.Range("I97").Value = 1143.1666
$ws_ARM.Range("K97").Value = 1143.1666
$ws_ARM.Range("M97").Value = -647.1666

# Row 102 (ARM)
$ws_ARM.Range("H102").Value = 1413.4615
$ws_ARM.Range("I102").Value = 1413.4615
$ws_ARM.Range("K102").Value = 1413.4615
$ws_ARM.Range("M102").Value = 208.5385000000001

# Row 136 (ARM)
$ws_ARM.Range("H136").Value = 6940.5713
$ws_ARM.Range("I136").Value = 6940.5713
$ws_ARM.Range("K136").Value = 20821.7139
$ws_ARM.Range("M136").Value = -18271.7139

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 22 (BSM)
$ws_BSM.Range("H22").Value = 557.5
$ws_BSM.Range("I22").Value = 557.5
$ws_BSM.Range("J22").Value = 0
$ws_BSM.Range("K22").Value = 557.5
$ws_BSM.Range("L22").Value = 0
$ws_BSM.Range("M22").ClearContents()
$ws_BSM.Range("N22").Value = -384.5

# Row 99 (BSM)
$ws_BSM.Range("H99").Value = 3489.125
$ws_BSM.Range("I99").Value = 3342
$ws_BSM.Range("K99").Value = 3342
$ws_BSM.Range("M99").Value = -1844

# Row 107 (BSM)
$ws_BSM.Range("H107").Value = 1942.625
$ws_BSM.Range("I107").Value = 1961.8695
$ws_BSM.Range("K107").Value = 1961.8695
$ws_BSM.Range("M107").Value = -41.86950000000002

# Row 135 (BSM)
$ws_BSM.Range("H135").Value = 70259
$ws_BSM.Range("J135").Value = 70259
$ws_BSM.Range("L135").Value = 70259
$ws_BSM.Range("N135").Value = -80399

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws_CRP.Range("H31").Value = 5748.9033
$ws_CRP.Range("I31").Value = 6120.2144
$ws_CRP.Range("K31").Value = 6120.2144
$ws_CRP.Range("M31").Value = -5825.2144

# Row 34 (CRP)
$ws_CRP.Range("H34").Value = 5748.9033
$ws_CRP.Range("I34").Value = 6120.2144
$ws_CRP.Range("K34").Value = 6120.2144
$ws_CRP.Range("M34").Value = -5918.2144

# Row 100 (CRP)
$ws_CRP.Range("H100").Value = 73000
$ws_CRP.Range("J100").Value = 73000
$ws_CRP.Range("L100").Value = 73000
$ws_CRP.Range("N100").Value = -75164

# Row 132 (CRP)
$ws_CRP.Range("H132").Value = 6197.364
$ws_CRP.Range("I132").Value = 5098.5293
$ws_CRP.Range("K132").Value = 15295.5879
$ws_CRP.Range("M132").Value = -12765.5879

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 64 (CUL)
$ws_CUL.Range("H64").Value = 8385.546
$ws_CUL.Range("I64").Value = 6428.8
$ws_CUL.Range("J64").Value = 10016.167
$ws_CUL.Range("K64").Value = 19286.4
$ws_CUL.Range("L64").Value = 30048.501
$ws_CUL.Range("M64").Value = -19016.4
$ws_CUL.Range("N64").Value = -30588.501

# Row 67 (CUL)
$ws_CUL.Range("H67").Value = 8385.546
$ws_CUL.Range("I67").Value = 6428.8
$ws_CUL.Range("J67").Value = 10016.167
$ws_CUL.Range("K67").Value = 19286.4
$ws_CUL.Range("L67").Value = 30048.501
$ws_CUL.Range("M67").Value = -18350.4
$ws_CUL.Range("N67").Value = -31920.501

# Row 109 (CUL)
$ws_CUL.Range("H109").Value = 2439.9092
$ws_CUL.Range("I109").Value = 1542.375
$ws_CUL.Range("K109").Value = 4627.125
$ws_CUL.Range("M109").Value = -3587.125

# Row 129 (CUL)
$ws_CUL.Range("H129").Value = 724.1429000000001
$ws_CUL.Range("I129").Value = 724.1429000000001
$ws_CUL.Range("J129").Value = 0
$ws_CUL.Range("K129").Value = 2172.4287
$ws_CUL.Range("L129").Value = 0
$ws_CUL.Range("M129").ClearContents()
$ws_CUL.Range("N129").Value = 2827.5713

# Row 141 (CUL)
$ws_CUL.Range("H141").Value = 5331.4287
$ws_CUL.Range("I141").Value = 4553.5
$ws_CUL.Range("K141").Value = 13660.5
$ws_CUL.Range("M141").Value = -8480.5

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 127 (GSM)
$ws_GSM.Range("H127").Value = 50561.4
$ws_GSM.Range("J127").Value = 50561.4
$ws_GSM.Range("L127").Value = 50561.4
$ws_GSM.Range("N127").Value = -60481.4

# Row 132 (GSM)
$ws_GSM.Range("H132").Value = 4244.4375
$ws_GSM.Range("I132").Value = 4316.2144
$ws_GSM.Range("J132").Value = 3742
$ws_GSM.Range("K132").Value = 12948.6432
$ws_GSM.Range("L132").Value = 11226
$ws_GSM.Range("M132").Value = -10418.6432
$ws_GSM.Range("N132").Value = -16286

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws_LTW.Range("H16").Value = 3837
$ws_LTW.Range("I16").Value = 2215.25
$ws_LTW.Range("J16").Value = 5999.3335
$ws_LTW.Range("K16").Value = 2215.25
$ws_LTW.Range("L16").Value = 5999.3335
$ws_LTW.Range("M16").Value = -2045.25
$ws_LTW.Range("N16").Value = -6339.3335

# Row 46 (LTW)
$ws_LTW.Range("H46").Value = 5953.222
$ws_LTW.Range("I46").Value = 6116.6
$ws_LTW.Range("J46").Value = 5749
$ws_LTW.Range("K46").Value = 6116.6
$ws_LTW.Range("L46").Value = 5749
$ws_LTW.Range("M46").Value = -5928.6
$ws_LTW.Range("N46").Value = -6125

# Row 93 (LTW)
$ws_LTW.Range("H93").Value = 2203.6667
$ws_LTW.Range("I93").Value = 1749.5
$ws_LTW.Range("J93").Value = 3112
$ws_LTW.Range("K93").Value = 1749.5
$ws_LTW.Range("L93").Value = 3112
$ws_LTW.Range("M93").Value = -501.5
$ws_LTW.Range("N93").Value = -5608

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws_WVR.Range("H81").Value = 23042.912
$ws_WVR.Range("J81").Value = 49438.555
$ws_WVR.Range("L81").Value = 98877.11
$ws_WVR.Range("N81").Value = -100999.11

# Row 84 (WVR)
$ws_WVR.Range("H84").Value = 23042.912
$ws_WVR.Range("J84").Value = 49438.555
$ws_WVR.Range("L84").Value = 494385.55
$ws_WVR.Range("N84").Value = -504993.55

# Row 96 (WVR)
$ws_WVR.Range("H96").Value = 1899.8
$ws_WVR.Range("I96").Value = 1899.75
$ws_WVR.Range("K96").Value = 1899.75
$ws_WVR.Range("M96").Value = -526.75

# Row 125 (WVR)
$ws_WVR.Range("H125").Value = 0
$ws_WVR.Range("J125").Value = 0
$ws_WVR.Range("L125").ClearContents()
$ws_WVR.Range("N125").Value = 0
